$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-word / reorder the "Butcher Shops" shared string used by B15 and
# append a trailing space (Excel will relocate it to the end of the
# shared-string table and shift every other appended-string index down).
$ws.Range("B15").Value = "(Butcher Shops and Canned Meat Factories) 1888-08-12 "

# Right-align the "Article" column (column B) data cells.
$ws.Range("B2:B33").HorizontalAlignment = -4152

# Move the active selection from C3 to B1.
$ws.Range("B1").Select()

"done"
